$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`nLodha Bhandup New Project - 2 & 3 BHK @ ₹ 2.29 Cr All inc`nnewproject-bhandup.co.in`nhttps://www.newproject-bhandup.co.in › lodha › new-launch`nLodha New Launch Luxurious 2 & 3 BHK apartments Starting Price from ₹ 2.29 Cr All Inc`nPrice List · Brochure · Request Site Visit · Configuration"
$ws.Range("B2").Value = "AddressofChoice Realty Pvt. Ltd"
$ws.Range("C2").Value = "India"

$ws.Range("A3").Value = "Sponsored`nLodha Bhandup | New Project Launch in Bhandup`nproptigermumbai.com`nhttps://www.proptigermumbai.com`nNew Launch Project — its Big Apartments with 2/3 BHK, its prime Location along the LBS...`nPricing & Floor Plan · Download Brochure · Runwal Group · Platinum Group · Adani Realty"
$ws.Range("B3").Value = "Locon Solutions Pvt. Ltd."
$ws.Range("C3").Value = "India"

$ws.Range("A4").Value = "Sponsored`nLödha Bhandup | Luxury 2, 2.5 & 3 BHK | Avail Pre Launch Offer`nl-bhandup.com`nhttps://www.l-bhandup.com › official`nExcellent Connectivity — Lödha on LBS Rd an upcoming residential project with the best of...`nBrochure & Floor Plans · Connect on Whatsapp · Price · View Project Highlights"
$ws.Range("B4").Value = "Rioga Premium Real Estate Advisory LLP"
$ws.Range("C4").Value = "India"

$ws.Range("A5").Value = "Sponsored`nNew Launch At Bhandup - New Launch At Bhandup, Mumbai`nlódháhomz.site`nhttps://www.lódháhomz.site › site-visit › enquire-now`nDownload Brochure — Its Big Apartments With 2/3 BHK, Its Prime Location Along The LBS..."
$ws.Range("B5").Value = "DIGITAL RUBIX"
$ws.Range("C5").Value = "India"

$ws.Range("A6").Value = "Sponsored`nLodha® coming soon to LBS`nlodhagroup.in`nhttps://www.lodhagroup.in`nLodha® LBS is Coming Soon — Live an exceptional lifestyle with forest living as Lodha comes to the prime LBS Marg. Lodha LBS along the Mulund-Bhandup-Kanjurmarg corridor ensures seamless connectivity."
$ws.Range("B6").Value = "Macrotech Developers Limited"
$ws.Range("C6").Value = "India"

$ws.Range("A7").Value = "Sponsored`nLodha Prelaunching Bhandup | 2 & 3 BHK Starting ₹2.29 Cr*`nprelaunch-projects.in`nhttps://www.prelaunch-projects.in`nGet EOI and Early Bird Benefits, Prime location at LBS Road Bhandup by Lodha. Lodha Bhandup 2 & 3 BHK Homes Starting Price...`nPrice Plans · Project Location · Site & Floor Plan · Virtual Site Visit · The Amenities"
$ws.Range("B7").Value = "PRELAUNCH REALTY PRIVATE LIMITED"
$ws.Range("C7").Value = "India"

$ws.Range("A8").Value = "Sponsored`nLodha New Launch Bhandup | 2/3 BHK Starts @ ₹ 2.29 Cr*`nbhandupnewlaunch.com`nhttps://www.bhandupnewlaunch.com › 2&3bhk › luxury_homes`nFree Pickup & Drop — Pre-book Lodha Bhandup at ₹1.08 Lacs* | Easy Access to Powai & R-City Mall | EOI Open Now! Modern Living at Lodha Bhandup | 10 mins to Eastern Express Hwy | Pre-book at ₹1.08 Lacs* Avail Special Offers."
$ws.Range("B8").Value = "Finwizz Holdings"
$ws.Range("C8").Value = "India"
